$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 96, shifting existing rows 96:149 down to 97:150
$ws.Rows.Item(96).Insert()

# Populate the new row 96 with the new record
$ws.Cells.Item(96, 1).Value = 11
$ws.Cells.Item(96, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(96, 3).Value = "Bíobío"
$ws.Cells.Item(96, 4).Value = 44518
$ws.Cells.Item(96, 5).Value = 8
$ws.Cells.Item(96, 6).Value = 100114013
$ws.Cells.Item(96, 7).Value = "Zanahoria"
$ws.Cells.Item(96, 8).Value = "Sin especificar"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 350
$ws.Cells.Item(96, 11).Value = 7000
$ws.Cells.Item(96, 12).Value = 7500
$ws.Cells.Item(96, 13).Value = 7286
$ws.Cells.Item(96, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(96, 15).Value = "Región Metropolitana"
$ws.Cells.Item(96, 16).Value = 364
$ws.Cells.Item(96, 17).Value = 20
$ws.Cells.Item(96, 18).Value = "Hortaliza"
